$wb = $excel.ActiveWorkbook

# Sheet 展览 (sheet1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Cells.Item(3, 6).Value = 34
$ws1.Cells.Item(5, 6).Value = 106
$ws1.Cells.Item(6, 2).Value = '2024-03-31'
$ws1.Cells.Item(6, 3).Value = '北京·GA06同人展'
$ws1.Cells.Item(6, 4).Value = '亦庄荣昌东街6号 北京亦创国际会展中心'
$ws1.Cells.Item(6, 5).Value = '2024.03.31 10:00-03.31 16:00'
$ws1.Cells.Item(6, 6).Value = 519
$ws1.Cells.Item(6, 7).Value = 85
$ws1.Cells.Item(6, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82620'
$ws1.Cells.Item(6, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/5VCyAnRb1710125054703.png'
$ws1.Cells.Item(7, 2).Value = '2024-04-04'
$ws1.Cells.Item(7, 3).Value = '北京·IDOx梦次元动漫游戏嘉年华3rd'
$ws1.Cells.Item(7, 4).Value = '北京展览馆 北京展览馆'
$ws1.Cells.Item(7, 5).Value = '2024.04.04 09:30-04.05 17:00'
$ws1.Cells.Item(7, 6).Value = 4907
$ws1.Cells.Item(7, 7).Value = 80
$ws1.Cells.Item(7, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=80825'
$ws1.Cells.Item(7, 9).Value = '//i1.hdslb.com/bfs/openplatform/202402/P1YCG3MT1708329896103.jpeg'
$ws1.Cells.Item(8, 6).Value = 4907
$ws1.Cells.Item(9, 3).Value = '北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 吴晛 专场活动'
$ws1.Cells.Item(9, 5).Value = '2024.04.04 10:30-04.04 13:30'
$ws1.Cells.Item(9, 6).Value = 114
$ws1.Cells.Item(9, 7).Value = 198
$ws1.Cells.Item(9, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82490'
$ws1.Cells.Item(9, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/ZVLr6IVF1709795299722.png'
$ws1.Cells.Item(10, 3).Value = '北京·IDOx梦次元动漫游戏嘉年华3rd同人创作大会'
$ws1.Cells.Item(10, 5).Value = '2024.04.04 09:30-04.05 17:00'
$ws1.Cells.Item(10, 6).Value = 131
$ws1.Cells.Item(10, 7).Value = 80
$ws1.Cells.Item(10, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82023'
$ws1.Cells.Item(10, 9).Value = '//i2.hdslb.com/bfs/openplatform/202402/DE1Xw4Ne1708668500346.png'
$ws1.Cells.Item(11, 3).Value = '北京·国乙同好嘉年华6th'
$ws1.Cells.Item(11, 4).Value = '北京国家会议中心 北京国家会议中心'
$ws1.Cells.Item(11, 5).Value = '2024.04.04 09:00-04.05 17:00'
$ws1.Cells.Item(11, 6).Value = 183
$ws1.Cells.Item(11, 7).Value = 85
$ws1.Cells.Item(11, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82372'
$ws1.Cells.Item(11, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/MyG450tb1709633846898.jpeg'
$ws1.Cells.Item(12, 3).Value = '北京·幻兽帕鲁only'
$ws1.Cells.Item(12, 4).Value = '北京展览馆 北京展览馆'
$ws1.Cells.Item(12, 5).Value = '2024.04.04 09:30-04.05 17:00'
$ws1.Cells.Item(12, 6).Value = 32
$ws1.Cells.Item(12, 7).Value = 80
$ws1.Cells.Item(12, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82549'
$ws1.Cells.Item(12, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/BbKUlDVR1709866539810.jpeg'
$ws1.Cells.Item(13, 3).Value = '北京·广播剧《你好撩人》专场活动'
$ws1.Cells.Item(13, 5).Value = '2024.04.04 11:35-04.04 14:50'
$ws1.Cells.Item(13, 6).Value = 208
$ws1.Cells.Item(13, 7).Value = 288
$ws1.Cells.Item(13, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82675'
$ws1.Cells.Item(13, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/oaz83RmQ1710152178775.png'
$ws1.Cells.Item(14, 3).Value = '北京·排球少年同好嘉年华'
$ws1.Cells.Item(14, 4).Value = '北京国家会议中心 北京国家会议中心'
$ws1.Cells.Item(14, 5).Value = '2024.04.04 09:30-04.05 17:00'
$ws1.Cells.Item(14, 6).Value = 141
$ws1.Cells.Item(14, 7).Value = 85
$ws1.Cells.Item(14, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82647'
$ws1.Cells.Item(14, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/p16AHCSz1710141324055.jpeg'
$ws1.Cells.Item(15, 3).Value = '北京·第15届IJOY漫展xCGF游戏节'
$ws1.Cells.Item(15, 5).Value = '2024.04.04 09:00-04.05 17:00'
$ws1.Cells.Item(15, 6).Value = 7945
$ws1.Cells.Item(15, 7).Value = 8.800000000000001
$ws1.Cells.Item(15, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81174'
$ws1.Cells.Item(15, 9).Value = '//i0.hdslb.com/bfs/openplatform/202401/EJejgoZa1705892035599.jpeg'
$ws1.Cells.Item(16, 6).Value = 7945
$ws1.Cells.Item(17, 3).Value = '北京·第四届花朝汉服节'
$ws1.Cells.Item(17, 4).Value = '南四环中路235号 世界花卉大观园'
$ws1.Cells.Item(17, 5).Value = '2024.04.04 10:00-04.06 17:00'
$ws1.Cells.Item(17, 6).Value = 261
$ws1.Cells.Item(17, 7).Value = 50
$ws1.Cells.Item(17, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82633'
$ws1.Cells.Item(17, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/wx3L2eSU1710137435377.jpeg'
$ws1.Cells.Item(18, 2).Value = '2024-04-05'
$ws1.Cells.Item(18, 3).Value = '北京·IDOx梦次元动漫游戏嘉年华3rd·配音演员 小N&小敢 专场活动'
$ws1.Cells.Item(18, 4).Value = '北京展览馆 北京展览馆'
$ws1.Cells.Item(18, 5).Value = '2024.04.05 10:30-04.05 13:45'
$ws1.Cells.Item(18, 6).Value = 134
$ws1.Cells.Item(18, 7).Value = 268
$ws1.Cells.Item(18, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82531'
$ws1.Cells.Item(18, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/S270ineo1709807616493.png'
$ws1.Cells.Item(19, 2).Value = '2024-04-06'
$ws1.Cells.Item(19, 3).Value = '北京·Yok运动番Only'
$ws1.Cells.Item(19, 4).Value = '宏福路53号 昆仑决世界搏击中心'
$ws1.Cells.Item(19, 5).Value = '2024.04.06 10:00-04.06 17:00'
$ws1.Cells.Item(19, 6).Value = 564
$ws1.Cells.Item(19, 7).Value = 78
$ws1.Cells.Item(19, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81595'
$ws1.Cells.Item(19, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/4lovHrUf1706759539872.jpeg'
$ws1.Cells.Item(20, 2).Value = '2024-04-13'
$ws1.Cells.Item(20, 3).Value = '北京·thebONE✖️GOJO超次元嘉年华12nd'
$ws1.Cells.Item(20, 4).Value = '小关路39号 北投购物公园'
$ws1.Cells.Item(20, 5).Value = '2024.04.13 10:00-04.14 17:00'
$ws1.Cells.Item(20, 6).Value = 2454
$ws1.Cells.Item(20, 7).Value = 63
$ws1.Cells.Item(20, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81708'
$ws1.Cells.Item(20, 9).Value = '//i1.hdslb.com/bfs/openplatform/202403/d216iQ1j1710843846207.jpeg'
$ws1.Cells.Item(21, 2).Value = '2024-04-19'
$ws1.Cells.Item(21, 3).Value = '北京·第22届中国国际模型博览会'
$ws1.Cells.Item(21, 4).Value = '北京展览馆 北京展览馆'
$ws1.Cells.Item(21, 5).Value = '2024.04.19 10:00-04.21 17:00'
$ws1.Cells.Item(21, 6).Value = 6306
$ws1.Cells.Item(21, 7).Value = 13.5
$ws1.Cells.Item(21, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=82425'
$ws1.Cells.Item(21, 9).Value = '//i2.hdslb.com/bfs/openplatform/202403/9nkCFSHm1709710888611.jpeg'
$ws1.Cells.Item(22, 2).Value = '2024-04-20'
$ws1.Cells.Item(22, 3).Value = '北京·QMQ动漫游戏嘉年华'
$ws1.Cells.Item(22, 4).Value = '小关路39号 北投购物公园'
$ws1.Cells.Item(22, 5).Value = '2024.04.20 10:00-04.21 17:00'
$ws1.Cells.Item(22, 6).Value = 2275
$ws1.Cells.Item(22, 7).Value = 70
$ws1.Cells.Item(22, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=81982'
$ws1.Cells.Item(22, 9).Value = '//i0.hdslb.com/bfs/openplatform/202402/lyPb1fLO1708569465126.jpeg'
$ws1.Cells.Item(23, 3).Value = '北京·亚力传感器走秀派对'
$ws1.Cells.Item(23, 4).Value = '旧鼓楼大街51号(鼓楼大街地铁站G东南口步行250米) MOONEE 暮霓'
$ws1.Cells.Item(23, 5).Value = '2024.04.20 19:00-04.21 02:00'
$ws1.Cells.Item(23, 6).Value = 2
$ws1.Cells.Item(23, 7).Value = 68
$ws1.Cells.Item(23, 8).Value = 'https://show.bilibili.com/platform/detail.html?id=83540'
$ws1.Cells.Item(23, 9).Value = '//i0.hdslb.com/bfs/openplatform/202403/ZDsD1X9t1711523212670.jpeg'
$ws1.Cells.Item(25, 6).Value = 2502
$ws1.Cells.Item(26, 6).Value = 14
$ws1.Cells.Item(27, 6).Value = 4
$ws1.Cells.Item(28, 6).Value = 6288
$ws1.Cells.Item(30, 6).Value = 52
$ws1.Cells.Item(31, 6).Value = 122
$ws1.Cells.Item(34, 6).Value = 6668
$ws1.Cells.Item(40, 6).Value = 103
$ws1.Cells.Item(43, 6).Value = 2489
$ws1.Cells.Item(47, 6).Value = 50
$ws1.Cells.Item(48, 6).Value = 471

# Sheet 本地生活 (sheet3)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Cells.Item(2, 6).Value = 1464

# Sheet 全部类型 (sheet4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Cells.Item(3, 6).Value = 1464
$ws4.Cells.Item(4, 6).Value = 34
$ws4.Cells.Item(7, 6).Value = 107
$ws4.Cells.Item(9, 6).Value = 519
$ws4.Cells.Item(10, 6).Value = 4907
$ws4.Cells.Item(11, 6).Value = 4907
$ws4.Cells.Item(13, 6).Value = 131
$ws4.Cells.Item(14, 6).Value = 183
$ws4.Cells.Item(17, 6).Value = 7945
$ws4.Cells.Item(18, 6).Value = 7945
$ws4.Cells.Item(21, 6).Value = 564
$ws4.Cells.Item(22, 6).Value = 2454
$ws4.Cells.Item(25, 6).Value = 6306
$ws4.Cells.Item(26, 6).Value = 2275
$ws4.Cells.Item(27, 6).Value = 2502
$ws4.Cells.Item(28, 6).Value = 14
$ws4.Cells.Item(30, 6).Value = 4
$ws4.Cells.Item(31, 6).Value = 6288
$ws4.Cells.Item(33, 6).Value = 52
$ws4.Cells.Item(34, 6).Value = 122
$ws4.Cells.Item(37, 6).Value = 6668
$ws4.Cells.Item(40, 6).Value = 103
$ws4.Cells.Item(43, 6).Value = 2489
$ws4.Cells.Item(46, 6).Value = 50
$ws4.Cells.Item(47, 6).Value = 471
